$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New model row labels were shuffled because a new "linear_svc" entry was
# added to the comparison (commit: "Incomplete svc model added").
# Row 8/9/10 now show naive bayes / linear_svc / perceptron respectively.
$ws.Range("A8").Value = "naive bayes"
$ws.Range("A9").Value = "linear_svc"
$ws.Range("A10").Value = "perceptron"

# Refreshed metric values (re-run of the model comparison with the new model).
$ws.Range("B2").Value = 0.95
$ws.Range("C2").Value = 0.56
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0.01
$ws.Range("F2").Value = 0.5607830844417628
$ws.Range("G2").Value = 0.5607830844417628

$ws.Range("B3").Value = 0.83
$ws.Range("C3").Value = 0.54
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 0.01
$ws.Range("F3").Value = 0.53677130971461
$ws.Range("G3").Value = 0.53677130971461

$ws.Range("B4").Value = 0.87
$ws.Range("C4").Value = 0.53
$ws.Range("D4").Value = 0.01
$ws.Range("E4").Value = 0.01
$ws.Range("F4").Value = 0.5336433037701915
$ws.Range("G4").Value = 0.5336433037701915

$ws.Range("B5").Value = 0.82
$ws.Range("C5").Value = 0.46
$ws.Range("D5").Value = 0.02
$ws.Range("E5").Value = 0.02
$ws.Range("F5").Value = 0.4586203688691701
$ws.Range("G5").Value = 0.4586203688691701

$ws.Range("B6").Value = 0.81
$ws.Range("C6").Value = 0.45
$ws.Range("D6").Value = 0.01
$ws.Range("E6").Value = 0.02
$ws.Range("F6").Value = 0.4454779218428099
$ws.Range("G6").Value = 0.4454779218428099

$ws.Range("B7").Value = 0.62
$ws.Range("C7").Value = 0.39
$ws.Range("D7").Value = 0.02
$ws.Range("E7").Value = 0.01
$ws.Range("F7").Value = 0.385767260489444
$ws.Range("G7").Value = 0.385767260489444

$ws.Range("B8").Value = 0.52
$ws.Range("C8").Value = 0.37
$ws.Range("D8").Value = 0.03
$ws.Range("E8").Value = 0.02
$ws.Range("F8").Value = 0.3702680024541976
$ws.Range("G8").Value = 0.3702680024541976

$ws.Range("B9").Value = 0.61
$ws.Range("C9").Value = 0.36
$ws.Range("D9").Value = 0.01
$ws.Range("E9").Value = 0.01
$ws.Range("F9").Value = 0.3563253812848866
$ws.Range("G9").Value = 0.3563253812848866

$ws.Range("B10").Value = 0.49
$ws.Range("C10").Value = 0.36
$ws.Range("D10").Value = 0.04
$ws.Range("E10").Value = 0.01
$ws.Range("F10").Value = 0.3563268787670251
$ws.Range("G10").Value = 0.3563268787670251
